$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns (default_count, default_value) right before the
#     existing "most_frequent_value" column (Q). Everything from Q onward
#     (most_frequent_value, memory_consumed_bytes, pattern_count, patterns)
#     shifts two columns to the right: Q->S, R->T, S->U, T->V.
$ws.Range("Q1:R1").EntireColumn.Insert()

# --- New header labels for the inserted columns.
$ws.Range("Q1").Value = "default_count"
$ws.Range("R1").Value = "default_value"

# --- Populate default_count / default_value for every data row.
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 17).Value = 0
    $ws.Cells.Item($row, 18).Value = "<Unspecified>"
}

# --- Refresh the most_frequent_value column (now column S) with its new
#     values for each attribute row.
$ws.Range("S2").Value = ""
$ws.Range("S3").Value = "Brislington"
$ws.Range("S4").Value = "The Avenue"
$ws.Range("S5").Value = "Pole"
$ws.Range("S6").Value = ""
$ws.Range("S7").Value = "No"
